$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67

$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 9.5
$ws.Range("AE6").Value = 15
$ws.Range("AN6").Value = 5.5
$ws.Range("AU6").Value = 8
